# Target: sheet "Action list" (the active sheet) - cell F4 is edited and
# becomes the new selection (was F6).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action list")

# F4 previously held "previous-sibling" (with a bold/applyFont style);
# it now holds "Test tj3 B" using the plain default style.
$ws.Range("F4").Value = "Test tj3 B"
$ws.Range("F4").Style = "Normal"

# The selection (active cell) moves from F6 to F4.
$ws.Range("F4").Select()
